$d = $word.ActiveDocument

# 1) Resize the syllabus table's three columns (widths are in twips in the
#    underlying XML; Word's COM Width property is expressed in points, so
#    divide by 20 to convert twips -> points).
$t = $d.Tables.Item(1)
$t.Columns.Item(1).Width = 953 / 20
$t.Columns.Item(2).Width = 1902 / 20
$t.Columns.Item(3).Width = 6735 / 20

# 2) Merge "Intro to " + "MATLAB" + ", variable types" into a single run of
#    text reading "Intro to MATLAB, variable types".
$d.Content.Find.Execute("Intro to MATLAB, variable types", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Intro to MATLAB, variable types", 2) | Out-Null

# 3) Replace the old lab-topics description with the updated biophysics text.
$d.Content.Find.Execute("modeling LIF, plotting place cells, or fitting simple GLMs", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Biophysics of single neurons and their modeling applications", 2) | Out-Null
